$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2 formulas
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3 formula
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15: shared formula for Area increments
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

$ws.Calculate()

# Match the final selection left by the author (G1 active cell, G1:K15 selected)
$ws.Range("G1:K15").Select() | Out-Null
